# Adds 7 new data rows (7-13) to the "Current Weather" sheet, mirroring the
# existing forecast-row layout, per commit "Modified the structure to add
# the new fields".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("A7").Value = "2017.05.29 03.41.21"
$ws.Range("B7").Value = 15.539999961853027
$ws.Range("C7").Value = 87.0
$ws.Range("D7").Value = 1014.0
$ws.Range("E7").Value = 1.5
$ws.Range("F7").Value = "2017-05-29T09:00:00"
$ws.Range("G7").Value = "2017-05-29T12:00:00"
$ws.Range("H7").Value = 14.079999923706055
$ws.Range("I7").Value = 995.1799926757812
$ws.Range("J7").Value = 90.0
$ws.Range("K7").Value = "'"
$ws.Range("K7").ClearFormats()
$ws.Range("L7").Value = 0.3100000023841858

# Row 8
$ws.Range("A8").Value = "2017.05.29 03.41.47"
$ws.Range("B8").Value = 15.539999961853027
$ws.Range("C8").Value = 87.0
$ws.Range("D8").Value = 1014.0
$ws.Range("E8").Value = 1.5
$ws.Range("F8").Value = "2017-05-29T09:00:00"
$ws.Range("G8").Value = "2017-05-29T12:00:00"
$ws.Range("H8").Value = 14.079999923706055
$ws.Range("I8").Value = 995.1799926757812
$ws.Range("J8").Value = 90.0
$ws.Range("K8").Value = "'"
$ws.Range("K8").ClearFormats()
$ws.Range("L8").Value = 0.3100000023841858

# Row 9
$ws.Range("A9").Value = "2017.05.29 03.42.21"
$ws.Range("B9").Value = 15.539999961853027
$ws.Range("C9").Value = 87.0
$ws.Range("D9").Value = 1014.0
$ws.Range("E9").Value = 1.5
$ws.Range("F9").Value = "2017-05-29T09:00:00"
$ws.Range("G9").Value = "2017-05-29T12:00:00"
$ws.Range("H9").Value = 14.079999923706055
$ws.Range("I9").Value = 995.1799926757812
$ws.Range("J9").Value = 90.0
$ws.Range("K9").Value = "'"
$ws.Range("K9").ClearFormats()
$ws.Range("L9").Value = 0.3100000023841858

# Row 10
$ws.Range("A10").Value = "2017.05.29 03.42.39"
$ws.Range("B10").Value = 15.539999961853027
$ws.Range("C10").Value = 87.0
$ws.Range("D10").Value = 1014.0
$ws.Range("E10").Value = 1.5
$ws.Range("F10").Value = "2017-05-29T09:00:00"
$ws.Range("G10").Value = "2017-05-29T12:00:00"
$ws.Range("H10").Value = 14.079999923706055
$ws.Range("I10").Value = 995.1799926757812
$ws.Range("J10").Value = 90.0
$ws.Range("K10").Value = "'"
$ws.Range("K10").ClearFormats()
$ws.Range("L10").Value = 0.3100000023841858

# Row 11
$ws.Range("A11").Value = "2017.05.29 03.43.42"
$ws.Range("B11").Value = 15.539999961853027
$ws.Range("C11").Value = 87.0
$ws.Range("D11").Value = 1014.0
$ws.Range("E11").Value = 1.5
$ws.Range("F11").Value = "2017-05-29T09:00:00"
$ws.Range("G11").Value = "2017-05-29T12:00:00"
$ws.Range("H11").Value = 14.079999923706055
$ws.Range("I11").Value = 995.1799926757812
$ws.Range("J11").Value = 90.0
$ws.Range("K11").Value = "'"
$ws.Range("K11").ClearFormats()
$ws.Range("L11").Value = 0.3100000023841858

# Row 12
$ws.Range("A12").Value = "2017.05.29 03.52.32"
$ws.Range("B12").Value = 15.539999961853027
$ws.Range("C12").Value = 87.0
$ws.Range("D12").Value = 1014.0
$ws.Range("E12").Value = 1.5
$ws.Range("F12").Value = "2017-05-29T09:00:00"
$ws.Range("G12").Value = "2017-05-29T12:00:00"
$ws.Range("H12").Value = 14.079999923706055
$ws.Range("I12").Value = 995.1799926757812
$ws.Range("J12").Value = 90.0
$ws.Range("K12").Value = "'"
$ws.Range("K12").ClearFormats()
$ws.Range("L12").Value = 0.3100000023841858

# Row 13 - structure shifted: K, L, M, N are blank (empty text) and the
# trailing numeric value lands in column O instead of L.
$ws.Range("A13").Value = "2017.05.29 03.58.20"
$ws.Range("B13").Value = 15.539999961853027
$ws.Range("C13").Value = 87.0
$ws.Range("D13").Value = 1014.0
$ws.Range("E13").Value = 1.5
$ws.Range("F13").Value = "2017-05-29T09:00:00"
$ws.Range("G13").Value = "2017-05-29T12:00:00"
$ws.Range("H13").Value = 14.079999923706055
$ws.Range("I13").Value = 995.1799926757812
$ws.Range("J13").Value = 90.0
$ws.Range("K13").Value = "'"
$ws.Range("K13").ClearFormats()
$ws.Range("L13").Value = "'"
$ws.Range("L13").ClearFormats()
$ws.Range("M13").Value = "'"
$ws.Range("M13").ClearFormats()
$ws.Range("N13").Value = "'"
$ws.Range("N13").ClearFormats()
$ws.Range("O13").Value = 0.3100000023841858

# The source workbook's declared dimension extends one column further
# (to P) than any populated cell - touch P13's formatting (a no-op value)
# so the sheet registers that column without writing visible content.
$ws.Range("P13").Font.Bold = $false
